$wb = $excel.ActiveWorkbook
$wsDados = $wb.Worksheets.Item("Dados")
$wsStats = $wb.Worksheets.Item("Estatisticas")

$aVals = @(
    0.003318612077168091,
    0.006978269609458462,
    0.01151923800666809,
    0.01647533303140074,
    0.02379583468112324,
    0.03064892338415876,
    0.03069441024871095,
    0.03688670165850495,
    0.0403176793981287,
    0.04053508883205179,
    0.1050936728213727,
    0.1060643747913208,
    0.1531513376774455,
    0.1996328961422279,
    0.2252565079552985,
    0.2641392101020513,
    0.2660901016169986,
    0.2852826224135205,
    0.2948130317003789,
    0.3000820616921689,
    0.3127722988103034,
    0.3209900740664595,
    0.3289369198204128,
    0.3890435267327689,
    0.5545456264128938,
    0.5845576760217531,
    0.6224487758558552,
    0.6652546731540427,
    0.9201614109964796,
    0.9252302812978552
)

$cVals = @(
    100062,
    100502,
    100298,
    100086,
    100207,
    100345,
    100573,
    100092,
    100456,
    100487,
    100512,
    100564,
    100485,
    100433,
    100221,
    100241,
    100330,
    100455,
    100467,
    100532,
    100004,
    100450,
    100347,
    100468,
    100538,
    100140,
    100460,
    100479,
    100462,
    100510
)

$bVals = @(
    "0.9998091982163766, 1.0000126972779593, 1.0000920464818637, 1.0000390595676432, 0.9999066828638551, 0.999712694971436, 1.0001200213352288, 1.0005797097282505, 1.0001441952238206, 0.9996422048248601, 1.001418516827778, 0.9991928994723775, 0.9986299300183125, 0.9995101057220779, 1.0001393854754341, 0.9998714028418239, 0.9995096098607855, 0.9995874354637853, 0.9999681679709743, 1.0003677484822913, 1.000257410524628, 1.000258699289715, 1.0002539356055435, 1.000221766976021, 1.0001171822820891, 0.9999127492735588, 0.999817706561212, 0.9998844456518291, 0.9999325000966167, 1.0000971589696859, 1.0001726501186121, 1.0001296319087183, 1.0000060476078934, 0.9998322425027588, 0.9995457007012626, 0.9990291692557287, 0.9994906823091751, 0.9998190724923317, 1.0001049580802999, 1.0001098361719307, 1.0002846904333833, 1.0000655749929706, 0.999967336569306, 1.000135690393142, 1.0002896450982228, 1.0005715479498958, 1.000761207110383, 1.0012898531183905, 1.0033317475767571, 1.0069933455159226",
    "1.0000222153725031, 0.9998915203222118, 1.0000133708854066, 1.0001340332690587, 0.9989650260244953, 0.9972776959343453, 0.9995042925018561, 1.0014690976919833, 1.0000494143407803, 0.9986800208766012, 1.0000210297560297, 1.0013303048789397, 0.999360291044547, 0.9998925902412855, 1.0003633510428607, 1.0003191833522371, 1.000427248392545, 0.999905107619002, 0.9993153841035205, 0.999875902971645, 1.0003329628213675, 0.9999505541102225, 1.0000430946258603, 1.0001162938775336, 1.0001052515749589, 1.0001412420965399, 0.9999035942718527, 0.9999338275242536, 0.9999444174026767, 1.0000363482300778, 1.0001405460860349, 1.0000258435769973, 0.9999134497228627, 0.9997439410489416, 0.9998102904614846, 0.9997879431454624, 0.9996593688097394, 0.9998884561958143, 1.0000529698509044, 1.0002384603861079, 1.0000349633249979, 0.9998499295372524, 1.0000234343097933, 1.0001458996293717, 1.0003074126171514, 1.000349109422251, 1.0005429552449647, 1.0010055168220322, 1.0022273715765653, 1.0045751412119501",
    "1.0002536368581756, 1.0032447004608787, 0.9977026351853443, 1.0004079705036544, 1.000398981430779, 1.0001701199693405, 1.0005077659332817, 1.001075508797672, 1.0003835695985466, 0.9998816019858883, 0.9999324519587866, 0.9999632610708817, 0.9999881366900074, 1.0000084366659436, 1.0000002955711282, 0.9999626398270117, 0.9998882216138405, 0.9999596093087282, 1.00002979330997, 1.0000829116814973, 1.000063561945785, 1.0000483245505056, 1.000040878174197, 1.0000524536236963, 1.0000340509733625, 0.9999639114012836, 0.9999634493446625, 0.9999827259539541, 0.9999870226861451, 1.000081163926809, 0.9999325281710897, 0.9999820840371996, 0.9999548436595158, 1.0000355405227328, 1.000176439048411, 1.0000143914794344, 0.9998462913198959, 0.9999347134412278, 0.9999850093547503, 1.0000121122469325, 0.9999765954908683, 0.9998867907008723, 0.9999352443915975, 1.000026158798927, 1.0001428795806628, 1.0000046383289765, 0.9998937404021866, 0.999715791355905, 0.9993389137469811, 0.998619878937119",
    "0.9999265430311315, 1.0000049757312004, 0.9999370557392878, 0.9998867512255143, 1.0001225635282345, 1.0005112136043481, 0.9997717588271597, 0.9991822292764242, 0.9999586614075062, 1.000842485696732, 1.000086358136424, 0.9993869861004477, 1.0007140399146157, 1.0001700240976703, 0.9997097037961943, 1.000472337059126, 1.0013795495492308, 1.0000385975823896, 0.998543700784399, 0.9998642797127476, 1.0012685286505947, 0.9992720591405695, 0.9996260846705683, 0.9997329665575383, 0.9999988386539604, 1.0004366516868866, 1.0011603227911965, 0.9989177889560162, 1.00013901566246, 1.0013649212008808, 1.0001478813937028, 0.9988504388562018, 1.0001599263105272, 1.001481466560638, 1.0000824533923707, 0.9988037470046043, 0.999710203042358, 1.0009315942206662, 1.0017473972621198, 1.0006951161309112, 0.9996748647519385, 0.9985018853492567, 0.9989786187535586, 0.9978888193750058, 0.9957723866508935, 0.9922689796399266, 0.9847858637749821, 0.9697800468490883, 0.9397855345285666, 0.8825303392787318",
    "1.0020789163673485, 0.999511237931841, 0.9998279166810605, 0.9999396885224648, 1.0000075438534497, 0.9999832886562252, 1.0000880687722735, 1.0002933334508184, 0.9999377959802671, 0.999485254328449, 0.9999308573278, 1.0005556939688791, 1.0016456624035173, 0.9997139922059977, 0.997615225897672, 1.0000386144909281, 1.0023818222071197, 0.99968277278076, 0.9967524759163835, 0.9990480647494724, 1.000924283523526, 1.003263363322311, 1.000204877987412, 0.9973603899273557, 0.9995193268950202, 1.0014073730399642, 1.0002163706757583, 0.9991781674131377, 0.999646751170616, 0.9997603967472527, 0.9998717404578148, 1.000088091282265, 1.0005224607647032, 1.0012934550241728, 1.0005813912577786, 1.0001370066075728, 1.0001517216915068, 0.9999700763581508, 0.9998997751919292, 0.9998797602806062, 0.9997777063774491, 0.9999131547692489, 0.999978733548969, 0.9999436173564937, 0.9999769654774748, 1.0004223339357723, 1.0010699715144478, 1.001490541943954, 1.002594036535906, 1.005357251187126",
    "1.0001874886559379, 1.0059804367907226, 1.0010974475063317, 0.9963431897034002, 1.0007016518373424, 0.9972214743592115, 0.9990133494724276, 1.0003103841210517, 1.0017556844765494, 1.0009009709421708, 1.0004247616558424, 1.0001432168326436, 0.9998794515583291, 0.999607252304201, 0.9990993476122959, 0.9996136100928481, 0.999829848949503, 0.9999544287971238, 1.0000158452398429, 1.0000393846023095, 1.0000849245837498, 1.000181309564799, 1.0003549652490191, 1.0001192116465296, 1.000005500653782, 0.9999349760477797, 0.9999654150611733, 1.0000079062688392, 1.0000298326976382, 0.999993746072536, 1.000018726781529, 1.000040011179113, 0.9999592359190251, 1.0000059598787197, 1.0000386843241145, 0.9999964730821878, 0.9999515363973347, 0.9999900672708257, 1.000014343218247, 1.0000236803921088, 1.0000126626392039, 1.0000200869112885, 1.000044980166582, 0.9999774834300419, 0.9999136071466362, 1.0000004033284242, 0.9998204907547296, 0.9995214669255513, 0.999281747979862, 0.998669642312096",
    "1.0000586161961142, 1.0004395090546407, 1.000166573390261, 0.9998365171967929, 1.000029274478206, 1.0000626099708085, 0.9996802270386971, 1.0011250896342851, 1.000056519465235, 0.9988850974167891, 0.9998560313854584, 1.0006648248693184, 0.9993184315600786, 0.9999895034429654, 1.0005176194647023, 0.9998047733097288, 0.9990762554537996, 1.001411045813969, 1.0045342334456946, 0.9964205348655953, 1.0009980653689161, 1.0008731504569799, 1.0002530606942706, 0.9998330067394897, 0.9998092290782424, 0.9997918441479667, 0.999769325611536, 0.9997744392159649, 0.9997981747903333, 0.9999224758598941, 0.9999919469163926, 1.0000652567980362, 0.9999956070546288, 1.0002690120039863, 1.000803559268585, 0.9999853294997507, 0.9989954770673823, 0.9998310670313973, 1.0002658931863913, 1.0009379311459248, 0.999825879712904, 0.9986112794842533, 0.9983720598619824, 0.9992557994838005, 0.999742831671866, 0.9966619672025275, 0.9944090477350929, 0.9894105091488826, 0.9785272160283344, 0.9571688305583543",
    "1.0026049743722414, 0.996461269465854, 1.0002915408612756, 1.0062574253978507, 1.0046657004862167, 0.9999658418260722, 1.0002001972428682, 1.000472262309553, 1.0003157532889904, 1.0001428918928028, 0.9999563236734753, 0.9996674678504263, 0.9991795980187559, 1.000057811643339, 1.0009922101800135, 1.0002985673390696, 0.9997550824751975, 0.9990730526276032, 0.9995785221287946, 0.9999604667843806, 1.000355479665387, 1.0009036801918831, 1.0003982532613898, 0.9999402641784932, 0.9993778023003538, 0.9997787314490106, 1.000103233050125, 1.0005095066635767, 1.0003280718669865, 1.0000979774706003, 0.9998930864494031, 1.0015064608758817, 1.0002469687124522, 0.999136689945785, 0.9996249281234821, 0.9999517537273723, 0.9998196762819365, 0.9995595126884285, 0.9997610939946874, 0.9998256433615069, 0.9998244531036062, 0.9997460967843506, 0.9993077704559005, 0.9984712113579157, 0.9969092645575854, 0.9948041490518168, 0.9899086186891731, 0.9799712910120792, 0.9599785925865395, 0.9211819036645166",
    "1.0000398172680902, 0.9999844839034897, 1.000112764096736, 1.0003013200257593, 0.9994925940303131, 0.9984464020428799, 1.0018170800931179, 0.9971343876071102, 0.9998622141368118, 1.0025232816803114, 1.000507633646882, 0.9987418717870334, 0.9995306689162738, 1.0000858045112055, 1.0006774767825057, 0.9975134830931944, 0.9930837302186732, 0.9985326680446422, 1.0032159646861762, 1.0008295168403243, 0.9988554251931875, 0.9992569489747763, 0.9992867934243295, 0.9999314948529096, 1.00055383202673, 1.0014640678243112, 1.000558748209057, 0.9999429877345022, 0.9992968889962621, 0.999860427985575, 1.000346946975035, 1.0001298400070837, 0.9999846892366365, 0.9998384589020501, 0.9996080399955084, 1.0002063743870686, 1.0000838013976994, 0.9999965582900883, 0.9999745781079669, 0.9999351751894396, 0.9999866384031284, 1.0000128638340156, 1.0000370643723664, 1.0000298378208752, 1.0000275594179329, 1.0000310239588899, 1.0001478644944384, 1.0003336172090953, 1.0006805738943216, 1.0013437239692073",
    "1.0071291331526713, 0.9980833790891005, 0.9991975397965395, 0.9998814423995414, 0.9997627310709094, 0.9995422963758913, 0.9994524375567528, 0.999063718410232, 1.0014355670084623, 1.004536202175033, 1.0008944982679104, 0.9977004863834487, 1.0001799582355158, 0.9996516427151699, 1.0002266251860448, 1.0008638135124424, 1.0005647464831324, 1.0004781981757276, 1.0006457634473298, 1.000200511051905, 0.9999027970143429, 0.9995803595196252, 0.9998178119940709, 0.9999897311251892, 1.0001603557505885, 0.9998724213425262, 0.9995174585195132, 0.9996589690608237, 0.9996891612638036, 0.9996356471890525, 0.9998520801038927, 1.000045024518036, 0.9995997518312464, 0.999887926647813, 1.0000640857682896, 0.9999711314092912, 0.9998566888312778, 1.0000133981648747, 1.0001281567608191, 1.0002878965526787, 0.9998875225136575, 1.0000125609271475, 1.000147859541716, 1.0002400482617415, 1.0004460014442427, 1.0007041094460334, 1.0013079667095468, 1.0031376958567384, 1.0065673617277455, 1.013336279643913",
    "0.9987118297655484, 0.9974578955603085, 1.0081190340541846, 0.9945492025960775, 0.9954318340733177, 1.0041345515080826, 0.9965355596015282, 0.9983618938080016, 0.9993597457710759, 0.9997783129953967, 1.0001519546016342, 0.9994039300497752, 0.9983298778007367, 0.9990511198023532, 0.9992853224166126, 1.0000907560913108, 1.0005110176850225, 1.0002929438071049, 1.0002148608708261, 1.0002423607639874, 1.0002123087955173, 1.0002673416249668, 1.0001731991776939, 1.0001372560614041, 1.0001047399511203, 1.0001418895785255, 1.0000308063522774, 0.9999707658662579, 0.9999576148107482, 1.000012621197873, 1.000008648736233, 1.0000042900655355, 0.9999878250981183, 1.0000153444847328, 1.000037241018674, 1.0000067723943464, 0.9999679985646615, 1.0000466635121001, 1.0001401216190449, 1.0003018795969385, 1.0002555410211564, 1.000102686685979, 0.9999956342063846, 1.0000157298014483, 1.0002203789488144, 1.0004952473559598, 1.0008835518253794, 1.001648373023743, 1.0033795646721768, 1.0068237473346624",
    "1.0013055364662469, 1.0064026674395803, 1.0115955478419885, 0.9956094692019019, 1.0027666038283807, 1.0002538490040869, 0.9978474464085699, 0.9996964294465888, 1.0013873154031865, 1.0037502495494375, 1.0013362105470123, 0.999587305556873, 0.9976192792757912, 0.9987493086457631, 0.9992414993258225, 1.0000385869734127, 1.0008421238743859, 1.0020609666324702, 1.0007198871766216, 0.9997451490274828, 1.000009096883432, 1.0002696474229675, 1.0006597090062395, 1.00017082829036, 0.9997615683056781, 0.9998470230192388, 0.9998892809010784, 0.9998654922877638, 0.9999748569464982, 0.9999881518722168, 0.9999882443899417, 1.0000053393316004, 1.0000118317267763, 1.0000107253335169, 1.0000045160613116, 0.9999929893401537, 1.0000025471017184, 1.0000074072482399, 1.0000059548790874, 0.9999995010091838, 0.9999881020702911, 0.9999650659417637, 1.0000048131214299, 1.000036820222297, 1.0000504562247896, 1.0000813992797177, 1.0001779232810826, 1.0003716954050006, 1.0007460275193811, 1.001475124511508",
    "1.000066294377897, 0.999432701831465, 0.9997163205536657, 0.9998320872525599, 0.9929398713539768, 0.9823844649653065, 0.9921614400022956, 0.997920722758828, 1.000512783207864, 1.0033232272127213, 1.001206085934646, 0.9996068638673008, 1.0004381027353257, 1.000943312817207, 1.0019020007316213, 0.9988672155429298, 0.9951016074380501, 1.005003911314111, 1.002351569455898, 1.0003993513881928, 0.9986462728017208, 0.99843655147877, 1.0007552899431202, 1.0034983322481736, 0.9999979032446301, 0.9964598847109342, 0.9994397176965135, 1.0020555008280991, 1.000898872897972, 1.0001759333057532, 1.0001462105297796, 1.0002362952614974, 1.0004632972919294, 1.0001392410397205, 0.9998646742012769, 0.9996766973982466, 0.999281828448789, 0.9996551882609119, 0.9998096714536372, 1.0001782253212765, 1.0006542906764309, 0.9998867567494564, 0.9999916255565883, 1.000061344248631, 0.9999330594239563, 0.999804895001213, 0.9996058822934504, 0.9993741484955024, 0.9988337993724389, 0.9975203153882402",
    "0.9997677084633059, 0.9996139291913719, 1.0003087596278553, 1.0014515270009954, 0.9961077666771009, 0.988770070958338, 1.0061113781085276, 1.008620843766567, 1.000422001571723, 0.9921655524229931, 1.0007529109669613, 1.0098211264316448, 1.003438146766508, 0.9987455431610622, 0.9972399492303701, 0.9999915409687319, 1.0027107025305257, 1.0005309750569729, 0.9986388218804746, 0.9991662944213014, 0.9993213684751612, 0.9994815233671666, 0.9994416363369045, 0.9990729312347214, 0.9993967281604565, 0.999344479753496, 0.9996323608694124, 0.9996138072153352, 0.9999904950116336, 1.0001817644936248, 1.0004192712157367, 1.0001624786116752, 1.0000588379967272, 1.0000088884635363, 0.999881185503588, 0.9995941280544225, 1.0001258059895681, 1.000029282134558, 1.0000013199788484, 0.999974574369908, 0.9999271170739221, 0.9999299853696857, 0.9999445411500237, 0.9999665494526522, 0.9996075441886939, 0.9989781232709279, 0.998886972039828, 0.9973410367891355, 0.9944343538892904, 0.9888547173579406",
    "1.0040370216977073, 1.0074971357699845, 0.9868454836818593, 0.9874490158213014, 0.9817785276021889, 0.990964813660146, 0.9954839018421586, 0.997841347834104, 0.9990329637404004, 0.9996254222713923, 0.9997193594828523, 0.9987487318139469, 0.9996449401308012, 1.0000969399169963, 1.0001040497714182, 1.0000603239653476, 1.000023222844117, 0.999887856237451, 0.9996755200573604, 0.9992909465817439, 0.9989302448901309, 0.9980549043239988, 0.999959319990468, 1.0018566731920207, 1.0001216297370166, 0.9994101874933762, 0.998351214954034, 1.0019134380200279, 1.006339636735084, 1.002276737250755, 0.9993706278814751, 0.9994990822234107, 0.9992910405131394, 1.0014397551601442, 1.0002694533003467, 0.9992374506875618, 0.9995388963063482, 0.9997215797923682, 0.9998439244015088, 0.9999227555976744, 0.9998844530212395, 0.9998812250328857, 1.0004385753847531, 1.0009560418018524, 1.001981680256128, 1.0041439205920395, 1.0084206147626271, 1.0164637663720306, 1.0330571526929393, 1.0677227732807655",
    "0.9905879458446973, 1.008021062108699, 0.9934961849093766, 1.00342523406428, 1.0001935721544624, 0.99541477607199, 0.9976521694989552, 0.998993750195673, 0.9994850852914835, 1.0001749477069184, 1.0005657773643233, 0.9954785813196341, 0.988080698069328, 0.9937590475501454, 0.9950875627575928, 0.9973937570581953, 0.9978289072921004, 0.9966296338024759, 0.9998174207955076, 1.0017023579679931, 0.9989078322520142, 0.9980012514685507, 0.9971874774489294, 1.001570012086434, 1.0063908788768354, 1.0027025632070206, 0.9997759068383193, 0.996728147047466, 0.9920038037372614, 0.9950570086272168, 0.99570030935748, 0.9976058117673581, 0.9989096688020236, 0.999828451563367, 1.0005745350205126, 1.0019966214504763, 1.0009484624525138, 1.0008864568507885, 1.0013884675287468, 0.9993813437508925, 0.9974456536213357, 0.9988488717489474, 0.9993415853966823, 0.999401789595055, 1.0006882919636604, 1.0019733812256686, 1.0024998044780977, 1.0036887140064996, 1.0066246881114342, 1.0157910204984792",
    "1.003656436959061, 1.018542588759574, 1.006162138723685, 0.9898452965055421, 0.9949252564826968, 0.9969662311227478, 0.9975688847580659, 1.0001049265273676, 1.002606547616906, 1.0018168429605983, 1.0006543690802299, 0.9998885695669727, 0.9991969376367055, 0.999167943894569, 0.9982684777265642, 0.9991908069311564, 0.9997448303368853, 1.0001307453765036, 1.0027635076928518, 1.0070846478090774, 0.9966523694581949, 0.9993805423125199, 0.9990575378069375, 0.9996292476718306, 1.0000510515465342, 1.0003121462126208, 1.0000590495513626, 0.9997136813275875, 0.9995519256465537, 0.999219116806578, 0.9986251600072823, 0.9974504032934098, 0.9957558284950592, 1.0043343807788987, 0.9981900799696105, 0.999375464322665, 1.000169422340626, 1.0003376005320954, 1.0003342403064193, 1.0003406738892167, 0.999806516808029, 0.9990472394721976, 1.0004940165620826, 1.002355133417106, 1.0029710650806336, 1.0057229986776985, 1.0120291151677507, 1.024200905964519, 1.049330714069298, 1.100716584613228",
    "1.004444173015775, 0.9989921955210691, 1.0221616451056907, 1.000941110607153, 1.0002915146870894, 0.9990153573301066, 0.9992746989630561, 0.9990371179184505, 0.9982332870419215, 0.9965768485893807, 1.0041706819795808, 1.002641612883382, 1.0024846356330825, 1.0000891185780503, 0.9976852909142133, 0.9987260454507667, 0.9992034831407268, 0.9994907374609538, 0.9994396224487975, 0.9988781300043394, 0.9995547309525976, 1.0000797291552965, 1.0006419926293628, 1.0003210225382186, 1.0000577086855171, 0.9999440772920166, 0.9997122141004755, 0.9997454220199766, 1.0001059131527557, 1.000175509367981, 1.000443035278427, 0.9999254068305565, 0.9998791600681605, 0.9999322326262258, 1.0001618708765387, 1.000863248057676, 1.0000660297807034, 0.9989995993771178, 0.9997035556549362, 1.0005604431283575, 1.0005177559204477, 1.0005122498650767, 1.0007123907682596, 1.0007152315227144, 1.0009213705838254, 1.0026277245112303, 1.0057584151092962, 1.0106893740298912, 1.021014976759974, 1.0423265662598662",
    "0.9738913551763232, 0.9737353044281563, 0.9896524760340665, 1.0004219034704411, 1.0020998765561773, 1.004621757105077, 1.0027195827675328, 1.0022089371799578, 0.997533262981217, 0.9987547545090728, 0.9992909661309992, 0.9994279285083498, 0.9994147667640392, 0.9993146823745538, 0.9989986078628385, 1.0001361267901405, 0.9994306581203798, 0.9996128176205272, 0.9994687034001989, 0.9998919769535397, 1.000072958453008, 1.0002802306500946, 1.0001112536912484, 0.9999606304821563, 0.9999900930988026, 1.0000846445843634, 0.9999876653032697, 0.9998543259602032, 0.9996266051154933, 0.9999811029049321, 1.0004545679450363, 0.99994847618263, 0.9996143924026338, 0.9999972744092093, 1.000189557924811, 1.0000413907586727, 0.9999949462290101, 1.0000372384952183, 1.0000265306610339, 0.9999934981816877, 0.9998516471194159, 1.000112923546464, 0.999978270573179, 1.0001238586495942, 1.0004902470250512, 1.000654523712937, 1.001201268997201, 1.0023402658581806, 1.0046343698437454, 1.0091534870407786",
    "0.9996677035452602, 1.0017328016627518, 0.9996088724618963, 0.9974940469991551, 1.0099033698914721, 1.0272736486810516, 1.010857205837213, 0.9999198536588916, 0.9988119194576679, 0.9970868085003426, 0.993851836660018, 0.9968796288407291, 0.9983620234198388, 1.001217166595517, 1.000124920398761, 0.9991909223876804, 0.9983835946563122, 0.9966715008232931, 0.9988168004442887, 1.0003031674296194, 1.0017310597858644, 0.9999450971553354, 0.9984168271151816, 0.9996675528300725, 1.0008760544029307, 0.99982150953551, 0.9988700187707141, 0.9995727760853559, 0.9999972756562217, 1.0006080932259245, 1.0002968098850122, 1.0006225911043507, 0.9998943218412191, 0.9990528515409385, 1.0000145802347347, 1.0008331957075398, 0.9999915380499571, 0.9994978783651997, 1.0001023938363913, 1.0003640544424817, 1.000854884317605, 1.0007171832719293, 0.9998432636217969, 0.9989485753385641, 0.9976159907695888, 0.995923879616557, 0.9921067418063776, 0.9838955412570464, 0.9677465618281957, 0.9360029885855312",
    "0.9998803609228534, 0.9999630317938408, 0.999624339864275, 0.9990697564940791, 1.0008572474163164, 1.0031484507422646, 0.9963462151659135, 0.9725001002725687, 0.986468856714932, 0.9936379679790951, 0.9981638866592059, 1.001598712343591, 0.999953648717449, 0.9984483045791235, 0.9961224914360304, 0.9981333603105084, 0.9991853264335163, 1.0010484353440665, 1.003440734050341, 1.001065747219285, 0.9992592773154417, 0.9976966898281656, 0.9949937317150116, 0.9993562910232632, 1.0033634668436557, 1.0008148227098792, 0.9986170286992687, 0.9957110879537292, 0.9982150494748391, 0.999827653330811, 1.0001982531056983, 1.0006492769133741, 1.0002928741163026, 1.0001049763509224, 1.000229197296288, 1.000459577317298, 1.0002415939119464, 1.000200975203394, 1.0002885297119968, 1.0008383005357506, 1.0018102010399161, 1.0004408856269547, 0.9992632803301279, 1.0008893731184518, 1.0029854682676813, 1.0061307071908252, 1.0123864929128878, 1.0234786460598158, 1.0469441876338876, 1.0967420653713127",
    "0.9798882902629581, 1.0047828374824244, 1.0029412987278232, 1.0021637391016842, 0.9997143732007182, 0.9972491333878033, 1.0004453972976945, 1.0068117051881933, 1.001103554306843, 0.9963456404842423, 0.9960413873744836, 0.9941771634300987, 0.9999471759313708, 1.00004835803195, 1.0008930999200594, 1.0024257137666017, 1.0053899558973451, 1.0109311269683998, 1.0056640701976969, 1.0033171467025672, 1.0015963502051406, 1.0008913250737572, 1.000821386367703, 0.9995354658764621, 0.9984094903539448, 0.9990940971563709, 0.9993273725680241, 1.0005328999689074, 1.0019647471644244, 1.000323016226731, 0.9989524511879111, 0.9972871481278334, 1.0005049314758292, 1.0037887658396363, 0.9994010541739862, 0.9945538435980878, 0.9981704626295934, 1.0009187857578814, 1.0043786058956068, 0.9997385290240943, 0.9998254668201111, 0.9998263178561102, 0.9998616429425128, 0.9999527373315369, 1.0000873960392105, 1.0003163903438537, 1.0007916107617947, 1.0018202185075642, 1.0045240811361282, 1.0095862811656644",
    "0.9999510049501245, 0.9903710254548012, 0.9979234182461293, 0.9901363913458031, 1.0195147953161248, 1.0085063268042649, 1.0018533682881403, 1.0018457719346505, 1.0020660493291862, 1.0031678949232155, 1.0015742430105925, 1.0007319532406633, 0.9998153429608609, 0.999053810531586, 0.9976958582421316, 0.9983517368770659, 1.0000743265512122, 1.0017180621580015, 1.0006614975217114, 1.0000619055744133, 0.9999655747126254, 1.0001822839106145, 0.9999016377224696, 0.9992593129915919, 0.9996666780425745, 0.9998672048181954, 1.0002625455110417, 1.0010309193025857, 1.0023630846083331, 0.9999574726942921, 0.9991930950506078, 0.9994554162783338, 0.9999517679559473, 1.0003920654694314, 1.000166733554661, 1.0001723468423522, 1.000416674324894, 1.000198588173702, 0.9999306034055987, 0.999934679978429, 0.999985925649144, 0.9989130969870883, 0.9988597123986076, 0.9982571598430868, 0.9963514533470952, 0.9927773878710044, 0.9855521220743642, 0.9711443301453819, 0.9426801876316773, 0.887993751061775",
    "0.9992757613584161, 1.00111657183418, 1.0000689226711497, 1.000211445663351, 0.9994179600690651, 0.9969778810635302, 1.0008378063079353, 1.0044975220904901, 0.9999061524077684, 0.9949100138338759, 0.9993886961083349, 1.0044723548000511, 0.9965646740511123, 0.9992086198474542, 1.0014498330958046, 0.9989544194806158, 0.9958764919175815, 1.0047274534167538, 1.0167048983500693, 1.000745569338523, 0.9848898766640827, 0.9948449892566822, 0.9993287421659415, 1.003493012706517, 0.9942923874918213, 1.0021682629727309, 1.011550859124128, 0.9994109139786729, 1.000169334166683, 1.0004706773345977, 1.000488668907994, 1.0007035886092088, 1.0016678330768976, 1.0014873946857832, 1.0018995839928673, 1.0031607320352673, 1.0003697538384928, 0.998657790465015, 0.9996746129065952, 1.000040710966746, 0.999954727739817, 0.9994250271607796, 0.9982353758554896, 1.0011577106330212, 1.0002559206102055, 0.9999070674208268, 1.0000496289196417, 1.000190057927208, 0.9989897554352378, 0.9983820473428945",
    "0.9999270764750308, 1.0012911346381332, 0.9998252258834468, 0.9982494619167858, 1.0100089457874297, 1.0268616427589965, 1.0056871407975725, 0.9871986811199865, 0.9990173763207142, 1.0102875531259354, 0.9994207346506154, 0.9881807974983509, 1.00557080730279, 1.0016547776183173, 0.9985590773670602, 0.9989993609472931, 0.9989181651809876, 0.9982800057200367, 0.9967665997896327, 0.9997614899507651, 0.9999154311288543, 1.000021522814592, 1.0000788003551142, 1.0001602598624748, 1.0003167667592014, 1.0002103346432243, 1.0002042756544263, 1.0001000978821744, 1.0000427575953124, 0.9999755836264236, 0.9998881213920461, 1.0005354406117197, 0.9997410055177107, 0.9988034429077584, 0.9998214959296858, 1.0007358309470697, 0.9995322767835394, 0.9997782465956607, 0.9999122339279802, 0.9999399022901047, 0.9999429444757674, 1.0000092101736828, 1.0000140109453106, 1.0000240603881474, 0.9999207481318794, 0.9997627347457855, 0.9994762678380884, 0.9989191678140185, 0.9978609998219692, 0.9957183033857412",
    "1.0047162824946676, 1.0283528324452071, 1.0004458828597658, 0.9923198386952632, 1.0116593862599328, 0.9994553706131857, 0.986875402764036, 0.9946604387634576, 0.9997237346281459, 1.0001778650343605, 1.0006452028794823, 1.00139034852038, 1.0028092515855411, 1.005888193671957, 1.0029515385910508, 1.0014717149903556, 1.000677502866232, 1.000170641354981, 0.9997067012735072, 0.999031674290959, 0.9978160099335008, 0.9988234197675716, 0.9992260170985676, 0.9995417224006272, 0.9996722927830222, 0.9996567947175479, 0.9999827248153045, 1.0003175422976793, 1.0002142872915172, 1.000229039347602, 1.0000862515120312, 0.9999637216071461, 1.0000608320906523, 1.0001601144668246, 1.0000568066136544, 0.999986453507848, 1.0002921920755294, 1.0007383115491757, 0.9983061140527872, 0.9994062276772459, 1.000204745634396, 0.9998890173163301, 0.9995412745260275, 0.998771268604521, 0.9974225140738298, 0.99588456308362, 0.9912380415164066, 0.9821425821720114, 0.9643130973026195, 0.9295936309035088",
    "1.0344441493951304, 1.0297869727244484, 0.9981941593146947, 1.0023737305275195, 1.0077036239569428, 1.005599890847949, 1.0002026104786255, 0.994879388085892, 0.9969017686843729, 0.9974649067997458, 0.9996280698628393, 1.001609149736418, 0.9979355777840726, 1.001272380411656, 1.005280440988801, 1.001648896566939, 0.9989564898659234, 0.9997698187012974, 0.9948578179446782, 0.9973580837620462, 0.9983429080808752, 0.9986625333766499, 0.9988949570729735, 0.9989599356456546, 0.999438866376604, 0.9995315778469542, 0.9993415362164706, 0.9998645165457914, 1.0004127225923927, 1.001307431442354, 1.0000374520703126, 0.9985901210856495, 0.9997981866360743, 1.00075323478514, 1.0002381439024761, 1.0000388591881402, 0.9989269791676747, 0.9994062298842917, 0.9995142857392236, 1.0004321765364754, 1.0014604350891327, 0.9992991546113678, 0.9994810096027822, 0.9995675836296376, 0.9992956711210603, 0.9985386337034006, 0.9968778465525754, 0.9939378651730154, 0.9880557651488385, 0.9750415280955155",
    "0.9999507854503583, 1.0017356515182267, 1.0003019487895155, 0.9990302665685875, 1.0161376956256791, 1.041528257940204, 1.016734874727392, 1.0004296531438173, 0.9986213933816197, 0.9963138016282639, 0.9973083516155341, 0.9968425235377675, 0.9990682042627836, 1.0010199277926175, 1.0033739066987848, 1.0027212652938566, 1.003420602590937, 0.9990063991955233, 0.9940173454241503, 0.9997877189130278, 1.0053151367830007, 1.0023172251431995, 1.0004879946686054, 0.9988709921314044, 0.9989177683894699, 0.9988735319254595, 0.9988066438300551, 0.9984719780777039, 0.9991907248693371, 0.9996998889135214, 1.0001091361762802, 1.0005559137713047, 1.0001008321608675, 1.0001860892416343, 1.0002418611235953, 0.9999918013002457, 0.9997878662624898, 0.9999609898958942, 0.9998979612863795, 1.0000830809368912, 1.0001155931818773, 1.0001847434816906, 1.0001102624569487, 0.9997105323752699, 1.0007149514207696, 0.99957922809149, 0.9982251907965212, 0.9993127438823702, 1.0000759261668037, 0.9982797451216797",
    "0.9619594005834968, 1.0078176667248961, 1.0025197558298091, 0.9984593767377475, 1.002230755029248, 1.005610547019428, 0.996834242800203, 0.9864049423828942, 1.000504241347355, 1.0147577305945075, 1.0102952512704924, 1.0033609326587791, 0.9981130478923349, 0.9988063867983431, 0.9989033086186845, 0.9984417105065994, 0.9971853592565445, 0.9945063173638967, 1.0012149695017414, 1.0006462016854052, 1.0004113179483234, 1.0003771736261056, 1.0004947590384903, 1.0007987427405296, 1.0015042019057696, 1.0005608055962718, 0.9999217854376536, 0.9997555713911849, 0.9994537163440756, 0.9998390156963769, 1.0001187901663242, 1.000546745202912, 1.0002085159064917, 0.9999954309596427, 0.9999196038467949, 0.9998116552284925, 1.0000209190246996, 1.0002760670799802, 1.000000596216561, 1.0000466803150976, 1.0001262287836012, 1.000151004681177, 1.0002182787366052, 1.0005153374764886, 1.001149481039605, 1.0023646772900325, 1.004511980918762, 1.0089698458228538, 1.0184290626283217, 1.0374993537025634",
    "0.9881281781210782, 1.0088083949535012, 0.9737416778104124, 0.9581111930921818, 0.9931882021531137, 0.9955827664801197, 0.9968231117657509, 0.9974433902295347, 0.9967637631866204, 0.9986363834195845, 0.9999077672990052, 1.0013855103344884, 1.0020173504544885, 0.995395235158987, 0.9974267296197283, 0.99876839202041, 0.9998546595440463, 1.0003387067423022, 0.9999331136268612, 0.9988334742596731, 0.9988267584422367, 0.9984694023533783, 0.999830102710702, 1.0011841647722266, 1.0031272256192412, 1.0002092993626923, 0.9975308087922985, 0.9995685412782401, 1.0014155349519887, 1.000701939509729, 1.0003117015575003, 0.9996772992143276, 0.9985943553463895, 0.9977982131070627, 0.9957812861066191, 0.9986989701919021, 0.9995582257023435, 0.9997176953128534, 0.9997147683528301, 0.9993849000684687, 1.0001458400222267, 1.0011125164408494, 0.9998310018586185, 0.9984774884149028, 0.9998667529110634, 1.0011865649571179, 1.0001487545164076, 0.9990995952540849, 0.9976138506862583, 0.9963313767618374"
)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $wsDados.Cells.Item($row, 1).Value = $aVals[$i]
    $wsDados.Cells.Item($row, 2).Value = $bVals[$i]
    $wsDados.Cells.Item($row, 3).Value = $cVals[$i]
}

$wsStats.Range("A2").Value = 0.003318612077168091
$wsStats.Range("B2").Value = 0.9252302812978552
$wsStats.Range("C2").Value = 0.2688239057002993
$wsStats.Range("D2").Value = 0.2593078926327553
$wsStats.Range("E2").Value = 100373.5333333333
